$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("work")

# Row 35 ("3" / Change sidebar presentation task row) was left half-filled;
# record the actual click-sidebar-interaction change task with its times.
$ws.Range("C35").Value = "Change click sidebar interaction"
$ws.Range("D35").Value = 0.84722222222222221   # 20:20 (was 18:20)
$ws.Range("E35").Value = 0.90972222222222221   # 21:50
$ws.Range("F35").Value = 1.5                   # hours

# New row 36: task "4" continuing that evening, starting 22:00
$ws.Range("B36").Value = "4"
$ws.Range("D36").Value = 0.91666666666666663   # 22:00

# Scroll the view down and leave the selection where the user left off.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("D37").Select()
